$d = $word.ActiveDocument
Write-Output ("paras before: " + $d.Paragraphs.Count)
$shp = $d.InlineShapes.Item(1)
$rng = $shp.Range
$insertPoint = $rng.Duplicate
$insertPoint.Collapse(0)
$insertPoint.InsertAfter("Z")
Write-Output ("paras mid1: " + $d.Paragraphs.Count)
$p94 = $d.Paragraphs.Item(94)
Write-Output "p94 text:[$($p94.Range.Text)] $($p94.Range.Start)-$($p94.Range.End)"

# Now delete just position 5281-5282 via a fresh Range object (not shp.Range)
$rng2 = $d.Range(5281, 5282)
Write-Output "rng2 text:[$($rng2.Text)]"
$rng2.Delete()
Write-Output ("paras mid2: " + $d.Paragraphs.Count)
$p94b = $d.Paragraphs.Item(94)
Write-Output "p94 text now:[$($p94b.Range.Text)] $($p94b.Range.Start)-$($p94b.Range.End)"
